# Natmi following Dr Hou advice
# Update the Col5a3-Sdc3 LR-pair sheet with recalculated NATMI statistics
# (Ligand/Receptor-expressing cell counts changed from 1 to 3, and all
# dependent expression / specificity metrics were recomputed accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used in this sheet (A..T):
#  A Sending cluster        K Receptor-expressing cells
#  B Ligand symbol          L Receptor detection rate
#  C Receptor symbol        M Receptor average expression value
#  D Target cluster         N Receptor total expression value
#  E Ligand-expressing cells  O Receptor derived specificity (avg)
#  F Ligand detection rate  P Receptor derived specificity (total)
#  G Ligand average expr    Q Edge average expression weight
#  H Ligand total expr      R Edge total expression weight
#  I Ligand derived spec(avg) S Edge average expression derived specificity
#  J Ligand derived spec(total) T Edge total expression derived specificity

$data = @{
  2  = @{ E=3; G=1.125230333333333; H=3.375691;      I=0.01107971898262852; J=0.01107971898262852; K=3; M=23.59622066666667; N=70.788662;        O=0.6996728317814862; P=0.6996728317814862; Q=26.55118324616022;  R=238.960649215442;   S=0.007752178355918786; T=0.007752178355918784 }
  3  = @{ E=3; G=1.125230333333333; H=3.375691;      I=0.01107971898262852; J=0.01107971898262852; K=3; M=7.778025666666667; N=23.334077;         O=0.2306332577891816; P=0.2306332577891816; Q=8.752070413578556;  R=78.768633722207;    S=0.002555351684352253; T=0.002555351684352253 }
  4  = @{ E=3; G=1.125230333333333; H=3.375691;      I=0.01107971898262852; J=0.01107971898262852; K=3; M=2.350402666666667; N=7.051208000000001;  O=0.06969391042933218; P=0.06969391042933218; Q=2.644744376080889; R=23.802699384728;    S=0.0007721889423574838; T=0.0007721889423574835 }
  5  = @{ E=3; G=85.94020333333333; H=257.82061;     I=0.8462207905669875; J=0.8462207905669875; K=3; M=23.59622066666667; N=70.788662;        O=0.6996728317814862; P=0.6996728317814862; Q=2027.864001991536; R=18250.77601792382; S=0.5920776968483721; T=0.5920776968483721 }
  6  = @{ E=3; G=85.94020333333333; H=257.82061;     I=0.8462207905669875; J=0.8462207905669875; K=3; M=7.778025666666667; N=23.334077;         O=0.2306332577891816; P=0.2306332577891816; Q=668.4451073252189; R=6016.00596592697;   S=0.1951666577374011; T=0.1951666577374011 }
  7  = @{ E=3; G=85.94020333333333; H=257.82061;     I=0.8462207905669875; J=0.8462207905669875; K=3; M=2.350402666666667; N=7.051208000000001;  O=0.06969391042933218; P=0.06969391042933218; Q=201.9940830885423; R=1817.94674779688;  S=0.05897643598121429; T=0.05897643598121429 }
  8  = @{ E=3; G=14.49222633333333; H=43.476679;     I=0.1426994904503839; J=0.1426994904503839; K=3; M=23.59622066666667; N=70.788662;        O=0.6996728317814862; P=0.6996728317814862; Q=341.9617705126109; R=3077.655934613498; S=0.09984295657719523; T=0.09984295657719523 }
  9  = @{ E=3; G=14.49222633333333; H=43.476679;     I=0.1426994904503839; J=0.1426994904503839; K=3; M=7.778025666666667; N=23.334077;         O=0.2306332577891816; P=0.2306332577891816; Q=112.7209083878092; R=1014.488175490283; S=0.03291124836742825; T=0.03291124836742825 }
  10 = @{ E=3; G=14.49222633333333; H=43.476679;     I=0.1426994904503839; J=0.1426994904503839; K=3; M=2.350402666666667; N=7.051208000000001;  O=0.06969391042933218; P=0.06969391042933218; Q=34.06256741980356; R=306.563106778232;  S=0.009945285505760397; T=0.009945285505760397 }
}

foreach ($row in $data.Keys) {
  $cols = $data[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = $cols[$col]
  }
}
